# STAI.xlsx update:
#  - fix typo "35. I feel inadequare" -> "35. I feel inadequate"
#  - add new "inverse" response-option value-label string and apply it to the
#    reverse-scored STAI items' Value Labels (column E)
#  - update selection to A15 (scrolled so row 9 is near the top)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$normalLabel  = "1=Not at all, 2=Somewhat, 3=Moderately so, 4=Very much so"
$reverseLabel = "4=Not at all, 3=Somewhat, 2=Moderately so, 1=Very much so"

# Rows (3-42) correspond to STAI questions 1-40. The reverse-scored items get
# the inverse response-option ordering in their Value Labels (column E); all
# other items keep the existing normal ordering.
$reverseRows = @(3,4,7,10,12,13,17,18,21,22,23,25,28,29,32,35,36,38,41)

for ($row = 3; $row -le 42; $row++) {
    if ($reverseRows -contains $row) {
        $ws.Range("E$row").Value = $reverseLabel
    } else {
        $ws.Range("E$row").Value = $normalLabel
    }
}

# Fix the question-35 typo.
$ws.Range("A37").Value = "35. I feel inadequate"

# Restore the workbook's on-open selection/scroll position.
$ws.Range("A15").Select()
$excel.ActiveWindow.ScrollRow = 9
